# Rename / re-word the header cells on the first sheet (tab "2025-06-08"):
#   C1: "agente" -> "agent"
#   D1: "numero"  -> "telefono"
#   I1: "ZIP CODE " -> "zip code"
# and move the active selection from K18 to I2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C1").Value = "agent"
$ws1.Range("D1").Value = "telefono"
$ws1.Range("I1").Value = "zip code"

$ws1.Range("I2").Select() | Out-Null
